# Generate Report for Handback
# Refresh the handback status timestamps / MT-vs-HT flags produced by the
# localization report generator.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview.Range("G3").Value = "2016-09-02 22:19:29"
$wsOverview.Range("G4").Value = "2016-09-02 22:19:29"

# --- zh-cn sheet ---
# Priority (column E): ht -> mt
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# Correspond Handoff Datetime (column H)
$wsZhCn.Range("H3").Value = "2016-09-02 22:19:25"
$wsZhCn.Range("H4").Value = "2016-09-02 22:19:25"

# Correspond Handback DateTime (column K)
$wsZhCn.Range("K3").Value = "2016-09-02 22:19:42"
$wsZhCn.Range("K4").Value = "2016-09-02 22:19:42"

# --- de-de sheet ---
# Priority (column E): ht -> mt
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# Correspond Handoff Datetime (column H)
$wsDeDe.Range("H3").Value = "2016-09-02 22:19:29"
$wsDeDe.Range("H4").Value = "2016-09-02 22:19:29"

# Correspond Handback DateTime (column K)
$wsDeDe.Range("K3").Value = "2016-09-02 22:19:49"
$wsDeDe.Range("K4").Value = "2016-09-02 22:19:49"
